$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object "object[,]" 36,16

$arr[0,0] = 2
$arr[0,1] = 1
$arr[0,2] = 65.652523
$arr[0,3] = 131.305046
$arr[0,4] = 0.07819491960606291
$arr[0,5] = 0.05521457487889056
$arr[0,6] = 2
$arr[0,7] = 1
$arr[0,8] = 0.7262925
$arr[0,9] = 1.452585
$arr[0,10] = 0.2952325527432663
$arr[0,11] = 0.2455349483544712
$arr[0,12] = 47.6829350609775
$arr[0,13] = 190.73174024391
$arr[0,14] = 0.02308568572685244
$arr[0,15] = 0.01355710779130248

$arr[1,0] = 2
$arr[1,1] = 1
$arr[1,2] = 65.652523
$arr[1,3] = 131.305046
$arr[1,4] = 0.07819491960606291
$arr[1,5] = 0.05521457487889056
$arr[1,6] = 3
$arr[1,7] = 1
$arr[1,8] = 0.4083613333333334
$arr[1,9] = 1.225084
$arr[1,10] = 0.1659958747772335
$arr[1,11] = 0.2070797486342548
$arr[1,12] = 26.80995182897734
$arr[1,13] = 160.859710973864
$arr[1,14] = 0.01298003408314386
$arr[1,15] = 0.0114338202868679

$arr[2,0] = 2
$arr[2,1] = 1
$arr[2,2] = 65.652523
$arr[2,3] = 131.305046
$arr[2,4] = 0.07819491960606291
$arr[2,5] = 0.05521457487889056
$arr[2,6] = 3
$arr[2,7] = 1
$arr[2,8] = 0.1071516666666667
$arr[2,9] = 0.321455
$arr[2,10] = 0.04355636342203116
$arr[2,11] = 0.05433653577813798
$arr[2,12] = 7.034777260321668
$arr[2,13] = 42.20866356193001
$arr[2,14] = 0.003405886336118186
$arr[2,15] = 0.003000168723381515

$arr[3,0] = 2
$arr[3,1] = 1
$arr[3,2] = 65.652523
$arr[3,3] = 131.305046
$arr[3,4] = 0.07819491960606291
$arr[3,5] = 0.05521457487889056
$arr[3,6] = 3
$arr[3,7] = 1
$arr[3,8] = 0.2805203333333333
$arr[3,9] = 0.841561
$arr[3,10] = 0.1140294497139816
$arr[3,11] = 0.1422516662860605
$arr[3,12] = 18.41686763613433
$arr[3,13] = 110.501205816806
$arr[3,14] = 0.008916523653108384
$arr[3,15] = 0.007854365279798638

$arr[4,0] = 2
$arr[4,1] = 1
$arr[4,2] = 65.652523
$arr[4,3] = 131.305046
$arr[4,4] = 0.07819491960606291
$arr[4,5] = 0.05521457487889056
$arr[4,6] = 3
$arr[4,7] = 1
$arr[4,8] = 0.1998293333333333
$arr[4,9] = 0.599488
$arr[4,10] = 0.08122915243236724
$arr[4,11] = 0.1013333162046457
$arr[4,12] = 13.11929990274133
$arr[4,13] = 78.71579941644801
$arr[4,14] = 0.006351707044117585
$arr[4,15] = 0.005595075975307703

$arr[5,0] = 2
$arr[5,1] = 1
$arr[5,2] = 65.652523
$arr[5,3] = 131.305046
$arr[5,4] = 0.07819491960606291
$arr[5,5] = 0.05521457487889056
$arr[5,6] = 2
$arr[5,7] = 1
$arr[5,8] = 0.7379140000000001
$arr[5,9] = 1.475828
$arr[5,10] = 0.29995660691112
$arr[5,11] = 0.2494637847424299
$arr[5,12] = 48.44591585702201
$arr[5,13] = 193.783663428088
$arr[5,14] = 0.02345508276272245
$arr[5,15] = 0.01377403682223233

$arr[6,0] = 3
$arr[6,1] = 1
$arr[6,2] = 276.9773456666667
$arr[6,3] = 830.932037
$arr[6,4] = 0.3298916825649746
$arr[6,5] = 0.3494120033757542
$arr[6,6] = 2
$arr[6,7] = 1
$arr[6,8] = 0.7262925
$arr[6,9] = 1.452585
$arr[6,10] = 0.2952325527432663
$arr[6,11] = 0.2455349483544712
$arr[6,12] = 201.1665688276075
$arr[6,13] = 1206.999412965645
$arr[6,14] = 0.09739476357242875
$arr[6,15] = 0.08579285820329813

$arr[7,0] = 3
$arr[7,1] = 1
$arr[7,2] = 276.9773456666667
$arr[7,3] = 830.932037
$arr[7,4] = 0.3298916825649746
$arr[7,5] = 0.3494120033757542
$arr[7,6] = 3
$arr[7,7] = 1
$arr[7,8] = 0.4083613333333334
$arr[7,9] = 1.225084
$arr[7,10] = 0.1659958747772335
$arr[7,11] = 0.2070797486342548
$arr[7,12] = 113.1068381795676
$arr[7,13] = 1017.961543616108
$arr[7,14] = 0.05476065842910641
$arr[7,15] = 0.07235614982884257

$arr[8,0] = 3
$arr[8,1] = 1
$arr[8,2] = 276.9773456666667
$arr[8,3] = 830.932037
$arr[8,4] = 0.3298916825649746
$arr[8,5] = 0.3494120033757542
$arr[8,6] = 3
$arr[8,7] = 1
$arr[8,8] = 0.1071516666666667
$arr[8,9] = 0.321455
$arr[8,10] = 0.04355636342203116
$arr[8,11] = 0.05433653577813798
$arr[8,12] = 29.67858421709278
$arr[8,13] = 267.1072579538351
$arr[8,14] = 0.01436888201570538
$arr[8,15] = 0.01898583782273754

$arr[9,0] = 3
$arr[9,1] = 1
$arr[9,2] = 276.9773456666667
$arr[9,3] = 830.932037
$arr[9,4] = 0.3298916825649746
$arr[9,5] = 0.3494120033757542
$arr[9,6] = 3
$arr[9,7] = 1
$arr[9,8] = 0.2805203333333333
$arr[9,9] = 0.841561
$arr[9,10] = 0.1140294497139816
$arr[9,11] = 0.1422516662860605
$arr[9,12] = 77.69777733219522
$arr[9,13] = 699.279995989757
$arr[9,14] = 0.03761736702810356
$arr[9,15] = 0.04970443970055162

$arr[10,0] = 3
$arr[10,1] = 1
$arr[10,2] = 276.9773456666667
$arr[10,3] = 830.932037
$arr[10,4] = 0.3298916825649746
$arr[10,5] = 0.3494120033757542
$arr[10,6] = 3
$arr[10,7] = 1
$arr[10,8] = 0.1998293333333333
$arr[10,9] = 0.599488
$arr[10,10] = 0.08122915243236724
$arr[10,11] = 0.1013333162046457
$arr[10,12] = 55.34819833300622
$arr[10,13] = 498.133784997056
$arr[10,14] = 0.02679682176924043
$arr[10,15] = 0.03540707702377402

$arr[11,0] = 3
$arr[11,1] = 1
$arr[11,2] = 276.9773456666667
$arr[11,3] = 830.932037
$arr[11,4] = 0.3298916825649746
$arr[11,5] = 0.3494120033757542
$arr[11,6] = 2
$arr[11,7] = 1
$arr[11,8] = 0.7379140000000001
$arr[11,9] = 1.475828
$arr[11,10] = 0.29995660691112
$arr[11,11] = 0.2494637847424299
$arr[11,12] = 204.3854610502727
$arr[11,13] = 1226.312766301636
$arr[11,14] = 0.09895318975039008
$arr[11,15] = 0.08716564079655034

$arr[12,0] = 3
$arr[12,1] = 1
$arr[12,2] = 219.2267506666667
$arr[12,3] = 657.680252
$arr[12,4] = 0.2611082919673688
$arr[12,5] = 0.2765585682093407
$arr[12,6] = 2
$arr[12,7] = 1
$arr[12,8] = 0.7262925
$arr[12,9] = 1.452585
$arr[12,10] = 0.2952325527432663
$arr[12,11] = 0.2455349483544712
$arr[12,12] = 159.22274480857
$arr[12,13] = 955.33646885142
$arr[12,14] = 0.0770876675799604
$arr[12,15] = 0.06790479376226696

$arr[13,0] = 3
$arr[13,1] = 1
$arr[13,2] = 219.2267506666667
$arr[13,3] = 657.680252
$arr[13,4] = 0.2611082919673688
$arr[13,5] = 0.2765585682093407
$arr[13,6] = 3
$arr[13,7] = 1
$arr[13,8] = 0.4083613333333334
$arr[13,9] = 1.225084
$arr[13,10] = 0.1659958747772335
$arr[13,11] = 0.2070797486342548
$arr[13,12] = 89.52372820457423
$arr[13,13] = 805.713553841168
$arr[13,14] = 0.0433428993367127
$arr[13,15] = 0.05726967878743968

$arr[14,0] = 3
$arr[14,1] = 1
$arr[14,2] = 219.2267506666667
$arr[14,3] = 657.680252
$arr[14,4] = 0.2611082919673688
$arr[14,5] = 0.2765585682093407
$arr[14,6] = 3
$arr[14,7] = 1
$arr[14,8] = 0.1071516666666667
$arr[14,9] = 0.321455
$arr[14,10] = 0.04355636342203116
$arr[14,11] = 0.05433653577813798
$arr[14,12] = 23.49051171185112
$arr[14,13] = 211.41460540666
$arr[14,14] = 0.01137292765743654
$arr[14,15] = 0.01502723453625745

$arr[15,0] = 3
$arr[15,1] = 1
$arr[15,2] = 219.2267506666667
$arr[15,3] = 657.680252
$arr[15,4] = 0.2611082919673688
$arr[15,5] = 0.2765585682093407
$arr[15,6] = 3
$arr[15,7] = 1
$arr[15,8] = 0.2805203333333333
$arr[15,9] = 0.841561
$arr[15,10] = 0.1140294497139816
$arr[15,11] = 0.1422516662860605
$arr[15,12] = 61.49756117259689
$arr[15,13] = 553.478050553372
$arr[15,14] = 0.02977403484879671
$arr[15,15] = 0.03934091715346582

$arr[16,0] = 3
$arr[16,1] = 1
$arr[16,2] = 219.2267506666667
$arr[16,3] = 657.680252
$arr[16,4] = 0.2611082919673688
$arr[16,5] = 0.2765585682093407
$arr[16,6] = 3
$arr[16,7] = 1
$arr[16,8] = 0.1998293333333333
$arr[16,9] = 0.599488
$arr[16,10] = 0.08122915243236724
$arr[16,11] = 0.1013333162046457
$arr[16,12] = 43.80793543455289
$arr[16,13] = 394.271418910976
$arr[16,14] = 0.02120960524957245
$arr[16,15] = 0.02802459684146119

$arr[17,0] = 3
$arr[17,1] = 1
$arr[17,2] = 219.2267506666667
$arr[17,3] = 657.680252
$arr[17,4] = 0.2611082919673688
$arr[17,5] = 0.2765585682093407
$arr[17,6] = 2
$arr[17,7] = 1
$arr[17,8] = 0.7379140000000001
$arr[17,9] = 1.475828
$arr[17,10] = 0.29995660691112
$arr[17,11] = 0.2494637847424299
$arr[17,12] = 161.7704884914427
$arr[17,13] = 970.6229309486561
$arr[17,14] = 0.07832115729489002
$arr[17,15] = 0.06899134712844959

$arr[18,0] = 3
$arr[18,1] = 1
$arr[18,2] = 136.2141876666667
$arr[18,3] = 408.642563
$arr[18,4] = 0.1622368336674611
$arr[18,5] = 0.1718366969192125
$arr[18,6] = 2
$arr[18,7] = 1
$arr[18,8] = 0.7262925
$arr[18,9] = 1.452585
$arr[18,10] = 0.2952325527432663
$arr[18,11] = 0.2455349483544712
$arr[18,12] = 98.93134289589251
$arr[18,13] = 593.588057375355
$arr[18,14] = 0.04789759455262924
$arr[18,15] = 0.04219191450346176

$arr[19,0] = 3
$arr[19,1] = 1
$arr[19,2] = 136.2141876666667
$arr[19,3] = 408.642563
$arr[19,4] = 0.1622368336674611
$arr[19,5] = 0.1718366969192125
$arr[19,6] = 3
$arr[19,7] = 1
$arr[19,8] = 0.4083613333333334
$arr[19,9] = 1.225084
$arr[19,10] = 0.1659958747772335
$arr[19,11] = 0.2070797486342548
$arr[19,12] = 55.6246072944769
$arr[19,13] = 500.621465650292
$arr[19,14] = 0.02693064512571875
$arr[19,15] = 0.03558390000417115

$arr[20,0] = 3
$arr[20,1] = 1
$arr[20,2] = 136.2141876666667
$arr[20,3] = 408.642563
$arr[20,4] = 0.1622368336674611
$arr[20,5] = 0.1718366969192125
$arr[20,6] = 3
$arr[20,7] = 1
$arr[20,8] = 0.1071516666666667
$arr[20,9] = 0.321455
$arr[20,10] = 0.04355636342203116
$arr[20,11] = 0.05433653577813798
$arr[20,12] = 14.59557723212945
$arr[20,13] = 131.360195089165
$arr[20,14] = 0.007066446487659558
$arr[20,15] = 0.009337010830147841

$arr[21,0] = 3
$arr[21,1] = 1
$arr[21,2] = 136.2141876666667
$arr[21,3] = 408.642563
$arr[21,4] = 0.1622368336674611
$arr[21,5] = 0.1718366969192125
$arr[21,6] = 3
$arr[21,7] = 1
$arr[21,8] = 0.2805203333333333
$arr[21,9] = 0.841561
$arr[21,10] = 0.1140294497139816
$arr[21,11] = 0.1422516662860605
$arr[21,12] = 38.21084932898255
$arr[21,13] = 343.897643960843
$arr[21,14] = 0.01849977686643936
$arr[21,15] = 0.02444405646585073

$arr[22,0] = 3
$arr[22,1] = 1
$arr[22,2] = 136.2141876666667
$arr[22,3] = 408.642563
$arr[22,4] = 0.1622368336674611
$arr[22,5] = 0.1718366969192125
$arr[22,6] = 3
$arr[22,7] = 1
$arr[22,8] = 0.1998293333333333
$arr[22,9] = 0.599488
$arr[22,10] = 0.08122915243236724
$arr[22,11] = 0.1013333162046457
$arr[22,12] = 27.21959031197156
$arr[22,13] = 244.976312807744
$arr[22,14] = 0.01317836049211881
$arr[22,15] = 0.01741278234447642

$arr[23,0] = 3
$arr[23,1] = 1
$arr[23,2] = 136.2141876666667
$arr[23,3] = 408.642563
$arr[23,4] = 0.1622368336674611
$arr[23,5] = 0.1718366969192125
$arr[23,6] = 2
$arr[23,7] = 1
$arr[23,8] = 0.7379140000000001
$arr[23,9] = 1.475828
$arr[23,10] = 0.29995660691112
$arr[23,11] = 0.2494637847424299
$arr[23,12] = 100.5143560778607
$arr[23,13] = 603.086136467164
$arr[23,14] = 0.0486640101428954
$arr[23,15] = 0.0428670327711046

$arr[24,0] = 3
$arr[24,1] = 1
$arr[24,2] = 66.466661
$arr[24,3] = 199.399983
$arr[24,4] = 0.07916459224847211
$arr[24,5] = 0.08384891234266049
$arr[24,6] = 2
$arr[24,7] = 1
$arr[24,8] = 0.7262925
$arr[24,9] = 1.452585
$arr[24,10] = 0.2952325527432663
$arr[24,11] = 0.2455349483544712
$arr[24,12] = 48.2742373843425
$arr[24,13] = 289.645424306055
$arr[24,14] = 0.02337196465639621
$arr[24,15] = 0.02058783836163373

$arr[25,0] = 3
$arr[25,1] = 1
$arr[25,2] = 66.466661
$arr[25,3] = 199.399983
$arr[25,4] = 0.07916459224847211
$arr[25,5] = 0.08384891234266049
$arr[25,6] = 3
$arr[25,7] = 1
$arr[25,8] = 0.4083613333333334
$arr[25,9] = 1.225084
$arr[25,10] = 0.1659958747772335
$arr[25,11] = 0.2070797486342548
$arr[25,12] = 27.14241430817467
$arr[25,13] = 244.281728773572
$arr[25,14] = 0.01314099574166813
$arr[25,15] = 0.0173634116911738

$arr[26,0] = 3
$arr[26,1] = 1
$arr[26,2] = 66.466661
$arr[26,3] = 199.399983
$arr[26,4] = 0.07916459224847211
$arr[26,5] = 0.08384891234266049
$arr[26,6] = 3
$arr[26,7] = 1
$arr[26,8] = 0.1071516666666667
$arr[26,9] = 0.321455
$arr[26,10] = 0.04355636342203116
$arr[26,11] = 0.05433653577813798
$arr[26,12] = 7.122013503918335
$arr[26,13] = 64.09812153526501
$arr[26,14] = 0.003448121750131362
$arr[26,15] = 0.004556059425464928

$arr[27,0] = 3
$arr[27,1] = 1
$arr[27,2] = 66.466661
$arr[27,3] = 199.399983
$arr[27,4] = 0.07916459224847211
$arr[27,5] = 0.08384891234266049
$arr[27,6] = 3
$arr[27,7] = 1
$arr[27,8] = 0.2805203333333333
$arr[27,9] = 0.841561
$arr[27,10] = 0.1140294497139816
$arr[27,11] = 0.1422516662860605
$arr[27,12] = 18.64524989927367
$arr[27,13] = 167.807249093463
$arr[27,14] = 0.009027094890925009
$arr[27,15] = 0.01192764749701728

$arr[28,0] = 3
$arr[28,1] = 1
$arr[28,2] = 66.466661
$arr[28,3] = 199.399983
$arr[28,4] = 0.07916459224847211
$arr[28,5] = 0.08384891234266049
$arr[28,6] = 3
$arr[28,7] = 1
$arr[28,8] = 0.1998293333333333
$arr[28,9] = 0.599488
$arr[28,10] = 0.08122915243236724
$arr[28,11] = 0.1013333162046457
$arr[28,12] = 13.28198855652267
$arr[28,13] = 119.537897008704
$arr[28,14] = 0.006430472730997339
$arr[28,15] = 0.008496688347834434

$arr[29,0] = 3
$arr[29,1] = 1
$arr[29,2] = 66.466661
$arr[29,3] = 199.399983
$arr[29,4] = 0.07916459224847211
$arr[29,5] = 0.08384891234266049
$arr[29,6] = 2
$arr[29,7] = 1
$arr[29,8] = 0.7379140000000001
$arr[29,9] = 1.475828
$arr[29,10] = 0.29995660691112
$arr[29,11] = 0.2494637847424299
$arr[29,12] = 49.046679685154
$arr[29,13] = 294.280078110924
$arr[29,14] = 0.02374594247835405
$arr[29,15] = 0.02091726701953633

$arr[30,0] = 2
$arr[30,1] = 1
$arr[30,2] = 75.063408
$arr[30,3] = 150.126816
$arr[30,4] = 0.08940367994566026
$arr[30,5] = 0.06312924427414179
$arr[30,6] = 2
$arr[30,7] = 1
$arr[30,8] = 0.7262925
$arr[30,9] = 1.452585
$arr[30,10] = 0.2952325527432663
$arr[30,11] = 0.2455349483544712
$arr[30,12] = 54.51799025483999
$arr[30,13] = 218.07196101936
$arr[30,14] = 0.02639487665499924
$arr[30,15] = 0.0155004357325082

$arr[31,0] = 2
$arr[31,1] = 1
$arr[31,2] = 75.063408
$arr[31,3] = 150.126816
$arr[31,4] = 0.08940367994566026
$arr[31,5] = 0.06312924427414179
$arr[31,6] = 3
$arr[31,7] = 1
$arr[31,8] = 0.4083613333333334
$arr[31,9] = 1.225084
$arr[31,10] = 0.1659958747772335
$arr[31,11] = 0.2070797486342548
$arr[31,12] = 30.652993375424
$arr[31,13] = 183.917960252544
$arr[31,14] = 0.01484064206088369
$arr[31,15] = 0.01307278803575975

$arr[32,0] = 2
$arr[32,1] = 1
$arr[32,2] = 75.063408
$arr[32,3] = 150.126816
$arr[32,4] = 0.08940367994566026
$arr[32,5] = 0.06312924427414179
$arr[32,6] = 3
$arr[32,7] = 1
$arr[32,8] = 0.1071516666666667
$arr[32,9] = 0.321455
$arr[32,10] = 0.04355636342203116
$arr[32,11] = 0.05433653577813798
$arr[32,12] = 8.04316927288
$arr[32,13] = 48.25901563728
$arr[32,14] = 0.003894099174980137
$arr[32,15] = 0.003430224440148718

$arr[33,0] = 2
$arr[33,1] = 1
$arr[33,2] = 75.063408
$arr[33,3] = 150.126816
$arr[33,4] = 0.08940367994566026
$arr[33,5] = 0.06312924427414179
$arr[33,6] = 3
$arr[33,7] = 1
$arr[33,8] = 0.2805203333333333
$arr[33,9] = 0.841561
$arr[33,10] = 0.1140294497139816
$arr[33,11] = 0.1422516662860605
$arr[33,12] = 21.056812233296
$arr[33,13] = 126.340873399776
$arr[33,14] = 0.01019465242660857
$arr[33,15] = 0.008980240189376413

$arr[34,0] = 2
$arr[34,1] = 1
$arr[34,2] = 75.063408
$arr[34,3] = 150.126816
$arr[34,4] = 0.08940367994566026
$arr[34,5] = 0.06312924427414179
$arr[34,6] = 3
$arr[34,7] = 1
$arr[34,8] = 0.1998293333333333
$arr[34,9] = 0.599488
$arr[34,10] = 0.08122915243236724
$arr[34,11] = 0.1013333162046457
$arr[34,12] = 14.999870778368
$arr[34,13] = 89.99922467020799
$arr[34,14] = 0.007262185146320611
$arr[34,15] = 0.006397095671791928

$arr[35,0] = 2
$arr[35,1] = 1
$arr[35,2] = 75.063408
$arr[35,3] = 150.126816
$arr[35,4] = 0.08940367994566026
$arr[35,5] = 0.06312924427414179
$arr[35,6] = 2
$arr[35,7] = 1
$arr[35,8] = 0.7379140000000001
$arr[35,9] = 1.475828
$arr[35,10] = 0.29995660691112
$arr[35,11] = 0.2494637847424299
$arr[35,12] = 55.390339650912
$arr[35,13] = 221.561358603648
$arr[35,14] = 0.026817224481868
$arr[35,15] = 0.01574846020455679

$ws.Range("E2:T37").Value = $arr
